# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.150.55'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '2.445.83'
$ws.Range('E3').Value = '  +0.53%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '580.43'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.31%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '143.29'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('D9').Value = '2.442.21'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('E11').Value = '  +2.39%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('E13').Value = '  -2.43%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '26.33'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').Value = '2.800.57'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').Value = '61.984.94'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '2.419.76'
$ws.Range('E18').Value = '  -0.52%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.83'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -3.22%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.17'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.19%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '328.37'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.81%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.09'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.49%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.99'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('E24').Value = '  +0.02%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '65.65'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.90%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.37'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +7.33%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '609.42'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0948'
$ws.Range('E29').Value = '  -5.38%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  -3.93%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.98'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  +1.02%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '4.89'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -3.78%  '
$ws.Range('E36').Value = '  +0.14%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.376'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  -4.90%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '149.23'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.85%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.31'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.40%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '18.35'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.63%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.73'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -2.13%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '42.56'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('E44').Value = '  +0.02%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.28%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '142.94'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('E47').Value = '  -2.85%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.604'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('E49').Value = '  -0.90%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '19.44'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -6.05%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +9.63%  '
